$d = $word.ActiveDocument

# Find the paragraph containing the attendance sentence, so the insertion is
# anchored to content rather than a hard-coded index.
$attendanceText = "Ursula, Jack, Nafees, Swaroop and Kanchan were in attendance"
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ($attendanceText + "`r")) {
        $targetIndex = $i
    }
}
if ($targetIndex -eq -1) {
    throw "Could not find the attendance paragraph"
}

# The paragraph right after it is the final, empty paragraph that precedes
# the section properties; insert the new signature-block paragraphs right
# before it (i.e. right after the attendance paragraph).
$finalPara = $d.Paragraphs.Item($targetIndex + 1)
$insertionPoint = $finalPara.Range
$insertionPoint.Collapse(1)

$xmlFragment = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>Please Print or Sign your signature here:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Ursula Mennear:  </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr><w:t>u mennear</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Cavolini" w:hAnsi="Cavolini" w:cs="Cavolini"/></w:rPr><w:tab/></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>Jack French:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>Nafees</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> Naushad Posharkar:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>Kanchan Dhansing Chavan:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>Kaiyuan Li:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:b/><w:bCs/></w:rPr><w:t>Swaroop Dattatraya Patil:</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertionPoint.InsertXML($xmlFragment)
Write-Host "Inserted signature block. New paragraph count:" $d.Paragraphs.Count
